# Updates the crypto price/volume table (Sheet1) to the latest scraped
# values, matching the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks like a plain number (e.g. "80.44")
# into a cell as literal text, instead of letting Excel auto-convert it
# to a floating point number (which would lose the exact original text,
# e.g. "80.439999999999998"). We briefly force a text number format,
# assign the value, then restore the cell's style so no visible
# formatting change is introduced.
function Set-TextValue {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "69.669.95"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "3.707.38"
$ws.Range("E3").Value = "  +0.75%  "
Set-TextValue $ws.Range("D5") "673.22"
$ws.Range("E5").Value = "  -1.29%  "
Set-TextValue $ws.Range("D6") "162.02"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("E12").Value = "  +1.44%  "
Set-TextValue $ws.Range("D13") "32.91"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "3.707.50"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "69.680.31"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("E17").Value = "  +2.31%  "
Set-TextValue $ws.Range("D18") "6.52"
$ws.Range("E18").Value = "  +2.07%  "
Set-TextValue $ws.Range("D19") "474.67"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("E21").Value = "  +0.86%  "
Set-TextValue $ws.Range("D22") "80.44"
$ws.Range("D23").Value = "3.855.29"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  +5.04%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +1.21%  "
Set-TextValue $ws.Range("D27") "9.18"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  +8.06%  "
Set-TextValue $ws.Range("D31") "2.01"
$ws.Range("E31").Value = "  +1.64%  "
Set-TextValue $ws.Range("D32") "6.59"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("E33").Value = "  +0.50%  "
Set-TextValue $ws.Range("D34") "0.999"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").Value = "3.698.01"
$ws.Range("E35").Value = "  +1.09%  "
Set-TextValue $ws.Range("D36") "8.54"
$ws.Range("E36").Value = "  +4.33%  "
$ws.Range("E37").Value = "  +1.56%  "
Set-TextValue $ws.Range("D39") "2.26"
$ws.Range("E39").Value = "  +1.98%  "
Set-TextValue $ws.Range("D40") "1.00"
$ws.Range("E40").Value = "  -0.05%  "
Set-TextValue $ws.Range("D41") "0.0915"
$ws.Range("E41").Value = "  +1.54%  "
Set-TextValue $ws.Range("D42") "173.95"
$ws.Range("E42").Value = "  +3.55%  "
$ws.Range("E43").Value = "  +0.16%  "
Set-TextValue $ws.Range("D44") "47.10"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("E51").Value = "  +0.68%  "

# Row 48/49 swap: SuiNetwork <-> InjectiveProtocol
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D48") "27.75"
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D49") "1.10"
$ws.Range("E49").Value = "  -0.52%  "
